$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.534.47"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "1.807.57"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.08"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.586"
$ws.Range("E6").Value = "  +4.83%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.87"
$ws.Range("E8").Value = "  +5.91%  "

$ws.Range("E9").Value = "  -0.41%  "

$ws.Range("E10").Value = "  -0.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0966"
$ws.Range("E11").Value = "  +1.54%  "

$ws.Range("D12").Value = "2.068.60"
$ws.Range("E12").Value = "  -0.58%  "

$ws.Range("E13").Value = "  +1.83%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.828.05"
$ws.Range("E14").Value = "  +0.49%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.653"
$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("E16").Value = "  +3.19%  "

$ws.Range("D17").Value = "34.498.34"
$ws.Range("E17").Value = "  -0.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.66"
$ws.Range("E18").Value = "  +0.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.48"
$ws.Range("E19").Value = "  -0.74%  "

$ws.Range("E20").Value = "  -1.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.63"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("E23").Value = "  -0.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.21"
$ws.Range("E24").Value = "  +6.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.66"
$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.97"
$ws.Range("E26").Value = "  +6.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.95"
$ws.Range("E27").Value = "  +1.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.119"
$ws.Range("E28").Value = "  +1.51%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.02"
$ws.Range("E30").Value = "  -0.74%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.85"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0530"
$ws.Range("E32").Value = "  -0.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.25"
$ws.Range("E33").Value = "  -0.20%  "

$ws.Range("E34").Value = "  -1.91%  "

$ws.Range("D35").Value = "1.396.09"
$ws.Range("E35").Value = "  -1.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.673"
$ws.Range("E36").Value = "  -0.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.46"
$ws.Range("E37").Value = "  -5.17%  "

$ws.Range("E38").Value = "  -0.41%  "

$ws.Range("E39").Value = "  -0.63%  "

$ws.Range("E40").Value = "  +1.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "83.23"
$ws.Range("E41").Value = "  -3.39%  "

$ws.Range("E42").Value = "  -0.67%  "

$ws.Range("E44").Value = "  +8.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.48"
$ws.Range("E45").Value = "  -2.40%  "

$ws.Range("E46").Value = "  -1.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0500"
$ws.Range("E47").Value = "  -5.11%  "

$ws.Range("D48").Value = "1.969.84"
$ws.Range("E48").Value = "  -0.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.62"
$ws.Range("E49").Value = "  -1.15%  "

$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("E51").Value = "  -2.82%  "

